# CancelacionAhorros.xlsx - "update entregable 1 y 2"
#
# - Rename header D1 from "razon" to "razoncierre"
# - Fill in the execution-result columns on row 2 (F:H): Estado, Transaccion, Fecha
# - Move the active selection to E11 (matches the saved sheet view state)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D1 header: razon -> razoncierre
$ws.Range("D1").Value = "razoncierre"

# Row 2 result columns
$ws.Range("F2").Value = "FAILED"

# G2 must stay a real (empty) text cell rather than being cleared outright.
# Writing a lone quote-prefix forces Excel to keep it as text; then resetting
# the style back to Normal drops the quote-prefix formatting while keeping
# the cell's empty string content.
$ws.Range("G2").Value = "'"
$ws.Range("G2").Style = "Normal"

$ws.Range("H2").Value = "26 jun. 2023, 18:04:05"

# Restore the selection/active cell recorded in the sheet view
$ws.Range("E11").Select()
